# Final data update for PR
# Populate the newly-logged workout times on the "Suivi" sheet. The
# "Point Hebdomadaire" and "Classement" sheets only contain formulas
# that reference "Suivi", so they recompute automatically once these
# values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suivi")

$updates = @{
    "BF2" = "17h31"
    "BG2" = "16h08"
    "BH2" = "18h26"
    "BI2" = "8h15"
    "BJ2" = "9h36"
    "BK2" = "7h09"
    "BL2" = "11h36"
    "BM2" = "10h58"
    "BN2" = "7h14"
    "BO2" = "17h58"
    "BP2" = "17h30"
    "BQ2" = "18h26"
    "BR2" = "6h40"
    "BS2" = "9h22"

    "BG3" = "21h07"
    "BH3" = "19h04"
    "BI3" = "20h42"
    "BN3" = "21h13"
    "BS3" = "20h44"

    "BF5" = "22h28"
    "BG5" = "21h35"
    "BH5" = "20h12"
    "BJ5" = "22h00"
    "BK5" = "22h18"

    "BI7" = "21h52"
    "BN7" = "23h11"
    "BO7" = "22h18"

    "BI8" = "7h40"

    "BF9" = "22h16"
    "BN9" = "22h13"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$excel.Calculate()
